# Update countries & provincias Spain
# Refresh the COVID-19 country stats table: a handful of rows swap which
# country name they show (two countries traded ranking positions in the
# shared-string table) and a batch of rows get refreshed case/death
# counters. Finally the "last updated" timestamp is bumped.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row-level data updates (country name swaps + refreshed statistics)
$rowUpdates = @(
    @{Row=4; B=4494742; C=61332; D=2174342; E=2168255; G=1070; H=152145}
    @{Row=5; B=2484649; C=41169; D=1721560; E=674455; G=955; H=88634}
    @{Row=6; B=1532135; C=49632; D=988770; E=509141}
    @{Row=17; B=267385; C=10284; D=136690; E=121621; G=297; H=9074}
    @{Row=23; B=173355; C=5939; E=95093; G=120; H=3179}
    @{Row=25; B=114994; C=397; D=100134; E=5948; G=11; H=8912}
    @{Row=50; B=41804; C=624; D=18764; E=22172; G=8; H=868}
    @{Row=56; A="Ghana"; B=34406; C=782; D=30621; E=3617; G=0; H=168}
    @{Row=57; A="Kirguistan"; B=33718; C=422; D=22296; E=10093; G=28; H=1329}
    @{Row=58; A="Japon"; B=30961; C=972; D=22811; E=7152; G=2; H=998}
    @{Row=59; A="Azerbaiyan"; B=30858; C=412; D=23873; E=6555; G=7; H=430}
    @{Row=70; A="Venezuela"; B=16571; C=583; D=10195; E=6225; G=5; H=151}
    @{Row=71; A="Costa Rica"; B=16344; C=503; D=3920; E=12299; G=10; H=125}
    @{Row=72; A="Chequia"; B=15799; C=283; D=11428; E=3997; G=1; H=374}
    @{Row=73; A="Costa de Marfil"; B=15713; C=58; D=10537; E=5078; G=2; H=98}
    @{Row=86; B=9150; C=18; E=143}
    @{Row=101; A="Paraguay"; B=4674; C=126; D=3039; E=1590; G=2; H=45}
    @{Row=102; A="Republica de Africa Central"; B=4599; D=1546; E=2994; H=59}
    @{Row=114; A="Montenegro"; B=2949; C=56; D=839; E=2065; H=45}
    @{Row=115; A="Mayotte"; B=2900; D=2672; E=190; H=38}
    @{Row=117; A="Zimbabue"; B=2817; C=113; D=604; E=2173; G=4; H=40}
    @{Row=118; A="Sri Lanka"; B=2810; C=5; D=2296; E=503; H=11}
    @{Row=137; B=1510; C=27; D=965; E=521}
    @{Row=146; B=1105; C=5; E=126}
    @{Row=163; A="Bahamas"; B=447; C=65; D=91; E=345; H=11}
    @{Row=164; A="Vietnam"; B=446; C=15; D=369; E=77; H=0}
    @{Row=165; A="Guyana"; B=389; D=181; E=188; H=20}
)

$colIndex = @{A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8}

foreach ($item in $rowUpdates) {
    $r = $item.Row
    foreach ($col in @("A","B","C","D","E","F","G","H")) {
        if ($item.ContainsKey($col)) {
            $ws.Cells.Item($r, $colIndex[$col]).Value = $item[$col]
        }
    }
}

# Updated timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 01:33"
